$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("A5:A6")
$dataRange.NumberFormat = "@"
$dataRange.Font.Bold = $true
$dataRange.HorizontalAlignment = -4108
$dataRange.VerticalAlignment = -4160
$dataRange.Borders.LineStyle = 1

$ws.Range("A5").Value = "08/02/2022"
$ws.Range("B5").Value = 3806.263471
$ws.Range("C5").Value = 11340.503841
$ws.Range("D5").Value = 0.55
$ws.Range("E5").Value = 10.65

$ws.Range("A6").Value = "09/02/2022"
$ws.Range("B6").Value = 3806.952118
$ws.Range("C6").Value = 11345.059095
$ws.Range("D6").Value = 0.54
$ws.Range("E6").Value = 10.65

$dataRange.NumberFormat = "General"
